$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing headers (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data for I2:J40 (I0 and IF values)
$iValues = @(5,8,7,6,9,9,2,6,9,6,6,10,6,7,9,6,6,7,6,9,9,6,9,9,5,9,7,5,7,6,6,7,8,6,5,5,6,5,3)
$jValues = @(6,8,8,6,9,9,3,7,9,6,6,10,7,7,9,6,6,8,7,9,9,8,9,9,6,9,7,6,8,6,6,7,8,6,5,5,6,5,3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

$wb.Save()
